$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column F previously used the plain default style; match it to column E's style ---
# (centered fill-only style -> text-number-format style used by the rest of the data block)
# Do this BEFORE writing any values so the "@" text format is already active and
# numeric-looking strings (e.g. "309750") are stored as text, not numbers.
$ws.Range("E1:E17").Copy() | Out-Null
$ws.Range("F1:F17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Column headers (row 1): standardballidx/oddballidx/bad_trials_* -> eyes_* ---
$ws.Range("C1").Value = "eyes_closed_start"
$ws.Range("D1").Value = "eyes_closed_stop"
$ws.Range("E1").Value = "eyes_open_start"
$ws.Range("F1").Value = "eyes_open_stop"

# --- Clear the old per-subject ball-task data (bad trial / idx lists) in C:F for every row ---
$ws.Range("C2:F17").ClearContents()

# --- New steady-state eyes-closed/open onset data, only populated for subject s1001 (row 2) ---
$ws.Range("C2").Value = "9500, 126750"
$ws.Range("D2").Value = "86500, 210500"
$ws.Range("E2").Value = "234000"
$ws.Range("F2").Value = "309750"

# --- Column widths: C/D keep their bestFit flag & grow to fit the longer header text;
#     E/F are brand-new bestFit columns sized for the new header/value text ---
$ws.Columns.Item(3).ColumnWidth = 19.998697916666668
$ws.Columns.Item(4).ColumnWidth = 19.666666666666668
$ws.Columns.Item(5).ColumnWidth = 13.830729166666666
$ws.Columns.Item(6).ColumnWidth = 13.498697916666666

# --- Selection moves to E8 ---
$ws.Range("E8").Select() | Out-Null
